$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected measurement values in column D ("2D") for ex_06, rows 2-33.
$ws.Range("D2").Value = 1000.239543
$ws.Range("D3").Value = 1025.052402
$ws.Range("D4").Value = 1055.262536
$ws.Range("D5").Value = 1006.16648
$ws.Range("D6").Value = 1014.45767
$ws.Range("D7").Value = 1036.4136189999999
$ws.Range("D8").Value = 1066.979335
$ws.Range("D9").Value = 1047.7962480000001
$ws.Range("D10").Value = 989.96896000000004
$ws.Range("D11").Value = 915.00840400000004
$ws.Range("D12").Value = 946.34419000000003
$ws.Range("D13").Value = 692.65493400000003
$ws.Range("D14").Value = 593.18312100000003
$ws.Range("D15").Value = 587.18239700000004
$ws.Range("D16").Value = 603.26649699999996
$ws.Range("D17").Value = 583.07719299999997
$ws.Range("D18").Value = 569.17149900000004
$ws.Range("D19").Value = 587.33351600000003
$ws.Range("D20").Value = 595.10094000000004
$ws.Range("D21").Value = 594.10118699999998
$ws.Range("D22").Value = 593.83093799999995
$ws.Range("D23").Value = 594.836232
$ws.Range("D24").Value = 581.74636099999998
$ws.Range("D25").Value = 589.20408299999997
$ws.Range("D26").Value = 597.84908700000005
$ws.Range("D27").Value = 577.07239100000004
$ws.Range("D28").Value = 588.29494399999999
$ws.Range("D29").Value = 580.92701599999998
$ws.Range("D30").Value = 583.76508100000001
$ws.Range("D31").Value = 584.70424100000002
$ws.Range("D32").Value = 585.98656100000005
$ws.Range("D33").Value = 564.50596800000005
